$wb = $excel.ActiveWorkbook

# ==========================================================================
# Sheet 1: "All Orders" -- a new order (#31) came in; insert a fresh row
# right under the header so the sheet stays newest-first, shifting every
# existing order down by one row.
# ==========================================================================
$ws1 = $wb.Worksheets.Item("All Orders")
$ws1.Rows.Item(2).Insert()

$ws1.Cells.Item(2, 1).Value = 31
$ws1.Cells.Item(2, 2).Value = "2026-01-27 14:09"
$ws1.Cells.Item(2, 3).Value = "Girija Lakade"
$ws1.Cells.Item(2, 4).Value = "A 1507"

# Phone number is purely numeric text -- force text so it isn't coerced to
# a number (matches the t="str" literal in every other Phone cell).
$phoneCell = $ws1.Cells.Item(2, 5)
$phoneCell.NumberFormat = "@"
$phoneCell.Value = "74996684"
$phoneCell.ClearFormats()

$ws1.Cells.Item(2, 6).Value = "Jawar Bhakari x4"
$ws1.Cells.Item(2, 7).Value = 80
$ws1.Cells.Item(2, 8).Value = "NEW"
$ws1.Cells.Item(2, 9).Value = "PENDING"

# Collection Date / Time / Notes / Cancel Reason / Feedback are still blank
# for this brand-new order, but stored as empty-string literals (not
# entirely absent cells) -- use the quote-prefix trick to force an empty
# text value, then strip the formatting it incidentally applies.
$blankRange = $ws1.Range("J2:N2")
$blankRange.Value = "'"
$blankRange.ClearFormats()

# ==========================================================================
# Sheet 2: "Daily Summary" -- add the 2026-01-27 roll-up row above
# 2026-01-26, shifting the rest of the history down by one row.
# ==========================================================================
$ws2 = $wb.Worksheets.Item("Daily Summary")
$ws2.Rows.Item(2).Insert()

# Date column holds text like "2026-01-26" elsewhere -- force text so COM
# doesn't coerce it into a date serial number.
$dateCell = $ws2.Cells.Item(2, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026-01-27"
$dateCell.ClearFormats()

$ws2.Cells.Item(2, 2).Value = 1
$ws2.Cells.Item(2, 3).Value = 0
$ws2.Cells.Item(2, 4).Value = 0
$ws2.Cells.Item(2, 5).Value = 80
$ws2.Cells.Item(2, 6).Value = 0
$ws2.Cells.Item(2, 7).Value = 80
